# Implement names, box-coloring, uptake threshold
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the sheet so row 2 is the top visible row (sheetView topLeftCell="A2"),
# without disturbing the current selection (stays on A1).
$excel.ActiveWindow.ScrollRow = 2

# Clear the "uptake" values that no longer apply (cells become blank, keeping style).
$ws.Range("C4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("C26").ClearContents()

# New uptake threshold value for row 11.
$ws.Range("C11").Value = 300
